$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"

$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "João Rodrigues-CAD"
$ws.Range("F3").Value = "Pedro Francisco-MTRM"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "[0, Andre Lucca-Acionamentos-2A, 0,"
$ws.Range("E4").Value = "João Rodrigues-CAD"

$ws.Range("C6").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "Josivaldo Ferreira-Programação"
$ws.Range("F6").Value = "-"

$ws.Range("C7").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("E7").Value = "Andre Barros-EAP"
$ws.Range("F7").Value = "Andre Barros-EAP"
